# Apply the edit described by the diff:
#  - Add new data to column D ("Made Changes" header + value 5) on the "Data" sheet
#  - Add a new, empty worksheet named "Rough" right after "Data"
#  - Keep "Data" as the active / selected sheet

$wb = $excel.ActiveWorkbook
$dataSheet = $wb.Worksheets.Item("Data")

# New column D: header in D1, value in D2
$dataSheet.Range("D1").Value = "Made Changes"
$dataSheet.Range("D2").Value = 5

# Auto-fit column D like Excel does when new data is entered
$dataSheet.Columns.Item(4).AutoFit() | Out-Null

# Insert a new blank worksheet named "Rough" immediately after "Data"
$roughSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $dataSheet)
$roughSheet.Name = "Rough"

# Keep "Data" as the active sheet
$dataSheet.Activate()
